# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: "Datos actualizados" timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 02:05"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1793263
$ws.Range("C4").Value = 24802
$ws.Range("D4").Value = 519381
$ws.Range("E4").Value = 1169343
$ws.Range("G4").Value = 1209
$ws.Range("H4").Value = 104539

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 468338
$ws.Range("C5").Value = 29526
$ws.Range("E5").Value = 247213
$ws.Range("G5").Value = 1180
$ws.Range("H5").Value = 27944

# --- Row 49: Panama ---
$ws.Range("B49").Value = 12531
$ws.Range("C49").Value = 400
$ws.Range("D49").Value = 7540
$ws.Range("E49").Value = 4665
$ws.Range("G49").Value = 6
$ws.Range("H49").Value = 326

# --- Row 190: San Vicente y las Granadinas ---
$ws.Range("B190").Value = 26
$ws.Range("C190").Value = 1
$ws.Range("E190").Value = 12

# --- Row 191: Gambia ---
$ws.Range("D191").Value = 20
$ws.Range("E191").Value = 4

# --- Rows 200/201: Santa Lucia <-> Belice swap (name + D/H values) ---
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# --- Rows 215/216: San Bartolome <-> Bonaire, San Eustaquio y Saba swap (name only, data identical) ---
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"
